# Apply the "Added Print_User_List and feedback funcs to excelfunc add 1 unittest" edit.
#
# Data changes:
#   - Users sheet:  B12 "snirben" -> "eladlp", D12 "sdf" -> "elad"  (unit-test user renamed)
#   - Cards sheet:  I3 "no" -> "yes"
#   - Games sheet:  new columns E ("date") / F ("feedback"); F2 gets sample feedback "it sucks"
#   - Games sheet becomes the active/selected sheet, with F1 selected.

$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B12").Value = "eladlp"
$wsUsers.Range("D12").Value = "elad"

$wsCards = $wb.Worksheets.Item("Cards")
$wsCards.Range("I3").Value = "yes"

$wsGames = $wb.Worksheets.Item("Games")
$wsGames.Range("E1").Value = "date"
$wsGames.Range("F1").Value = "feedback"
$wsGames.Range("F2").Value = "it sucks"

$wsGames.Columns.Item(5).ColumnWidth = 16
$wsGames.Columns.Item(6).ColumnWidth = 21.333333333333332

$wsGames.Activate() | Out-Null
$wsGames.Range("F1").Select() | Out-Null
